$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.52%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'37.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.46%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.90%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07857"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.38%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.183"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-8.29%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'8.004"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'4.004"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.59%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9100"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.07%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1876"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.12%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09214"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.93%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.21%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03517"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.09%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.09942"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.50%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001491"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.35%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005661"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.82%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.479"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.51%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'-1.92%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'2.84%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1314"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.49%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.543"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.33%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2241"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.54%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04635"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.22%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001227"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.82%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004444"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.35%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001300"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0004748"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'39.94%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.01746"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.34%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04714"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.66%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007890"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.85%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1391"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.39%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.007661"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.91%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.50%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.01022"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'11.24%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006060"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.05%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'8.668"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'183.00%"
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'34.83%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.12%"
$ws.Range("E51").Style = "Normal"
